$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the J2 RF connector part-number references (row 7) to the new
# 142-0711-301 connector, replacing the old SMB_131-3711-301 /
# CONN_131-3711-301 references. A leading apostrophe is used so the
# cells keep being stored/entered as explicit text (preserving their
# existing "quote prefix" text style) instead of Excel re-evaluating them.
$ws.Range("B7").Formula = "'142-0711-301"
$ws.Range("D7").Formula = "'CONN_142-0711-301"
$ws.Range("E7").Formula = "'CONN_142-0711-301"

# Leave the active selection on B10, matching the saved workbook view state.
$ws.Range("B10").Select()
